$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(160, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(160, 2).Value = 'PASSED'
$ws.Cells.Item(160, 3).Value = 'chrome'
$ws.Cells.Item(160, 4).Formula = '="02.11.22"'
$ws.Cells.Item(160, 4).Copy()
$ws.Cells.Item(160, 4).PasteSpecial(-4163)

$ws.Cells.Item(161, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(161, 2).Value = 'PASSED'
$ws.Cells.Item(161, 3).Value = 'chrome'
$ws.Cells.Item(161, 4).Formula = '="02.11.22"'
$ws.Cells.Item(161, 4).Copy()
$ws.Cells.Item(161, 4).PasteSpecial(-4163)

$ws.Cells.Item(162, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(162, 2).Value = 'PASSED'
$ws.Cells.Item(162, 3).Value = 'chrome'
$ws.Cells.Item(162, 4).Formula = '="02.11.22"'
$ws.Cells.Item(162, 4).Copy()
$ws.Cells.Item(162, 4).PasteSpecial(-4163)

$ws.Cells.Item(163, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(163, 2).Value = 'PASSED'
$ws.Cells.Item(163, 3).Value = 'chrome'
$ws.Cells.Item(163, 4).Formula = '="02.11.22"'
$ws.Cells.Item(163, 4).Copy()
$ws.Cells.Item(163, 4).PasteSpecial(-4163)

$ws.Cells.Item(164, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(164, 2).Value = 'PASSED'
$ws.Cells.Item(164, 3).Value = 'chrome'
$ws.Cells.Item(164, 4).Formula = '="02.11.22"'
$ws.Cells.Item(164, 4).Copy()
$ws.Cells.Item(164, 4).PasteSpecial(-4163)

$ws.Cells.Item(165, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(165, 2).Value = 'PASSED'
$ws.Cells.Item(165, 3).Value = 'chrome'
$ws.Cells.Item(165, 4).Formula = '="02.11.22"'
$ws.Cells.Item(165, 4).Copy()
$ws.Cells.Item(165, 4).PasteSpecial(-4163)

$ws.Cells.Item(166, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(166, 2).Value = 'PASSED'
$ws.Cells.Item(166, 3).Value = 'chrome'
$ws.Cells.Item(166, 4).Formula = '="02.11.22"'
$ws.Cells.Item(166, 4).Copy()
$ws.Cells.Item(166, 4).PasteSpecial(-4163)

$ws.Cells.Item(167, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(167, 2).Value = 'PASSED'
$ws.Cells.Item(167, 3).Value = 'chrome'
$ws.Cells.Item(167, 4).Formula = '="02.11.22"'
$ws.Cells.Item(167, 4).Copy()
$ws.Cells.Item(167, 4).PasteSpecial(-4163)

$ws.Cells.Item(168, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(168, 2).Value = 'PASSED'
$ws.Cells.Item(168, 3).Value = 'chrome'
$ws.Cells.Item(168, 4).Formula = '="02.11.22"'
$ws.Cells.Item(168, 4).Copy()
$ws.Cells.Item(168, 4).PasteSpecial(-4163)

$ws.Cells.Item(169, 1).Value = 'Login with valid username and password'
$ws.Cells.Item(169, 2).Value = 'PASSED'
$ws.Cells.Item(169, 3).Value = 'chrome'
$ws.Cells.Item(169, 4).Formula = '="02.11.22"'
$ws.Cells.Item(169, 4).Copy()
$ws.Cells.Item(169, 4).PasteSpecial(-4163)

$ws.Cells.Item(170, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(170, 2).Value = 'PASSED'
$ws.Cells.Item(170, 3).Value = 'chrome'
$ws.Cells.Item(170, 4).Formula = '="02.11.22"'
$ws.Cells.Item(170, 4).Copy()
$ws.Cells.Item(170, 4).PasteSpecial(-4163)

$ws.Cells.Item(171, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(171, 2).Value = 'PASSED'
$ws.Cells.Item(171, 3).Value = 'chrome'
$ws.Cells.Item(171, 4).Formula = '="02.11.22"'
$ws.Cells.Item(171, 4).Copy()
$ws.Cells.Item(171, 4).PasteSpecial(-4163)

$ws.Cells.Item(172, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(172, 2).Value = 'PASSED'
$ws.Cells.Item(172, 3).Value = 'chrome'
$ws.Cells.Item(172, 4).Formula = '="02.11.22"'
$ws.Cells.Item(172, 4).Copy()
$ws.Cells.Item(172, 4).PasteSpecial(-4163)

$ws.Cells.Item(173, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(173, 2).Value = 'PASSED'
$ws.Cells.Item(173, 3).Value = 'chrome'
$ws.Cells.Item(173, 4).Formula = '="02.11.22"'
$ws.Cells.Item(173, 4).Copy()
$ws.Cells.Item(173, 4).PasteSpecial(-4163)

$ws.Cells.Item(174, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(174, 2).Value = 'PASSED'
$ws.Cells.Item(174, 3).Value = 'chrome'
$ws.Cells.Item(174, 4).Formula = '="02.11.22"'
$ws.Cells.Item(174, 4).Copy()
$ws.Cells.Item(174, 4).PasteSpecial(-4163)

$ws.Cells.Item(175, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(175, 2).Value = 'PASSED'
$ws.Cells.Item(175, 3).Value = 'chrome'
$ws.Cells.Item(175, 4).Formula = '="02.11.22"'
$ws.Cells.Item(175, 4).Copy()
$ws.Cells.Item(175, 4).PasteSpecial(-4163)

$ws.Cells.Item(176, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(176, 2).Value = 'FAILED'
$ws.Cells.Item(176, 3).Value = 'chrome'
$ws.Cells.Item(176, 4).Formula = '="02.11.22"'
$ws.Cells.Item(176, 4).Copy()
$ws.Cells.Item(176, 4).PasteSpecial(-4163)

$ws.Cells.Item(177, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(177, 2).Value = 'FAILED'
$ws.Cells.Item(177, 3).Value = 'chrome'
$ws.Cells.Item(177, 4).Formula = '="02.11.22"'
$ws.Cells.Item(177, 4).Copy()
$ws.Cells.Item(177, 4).PasteSpecial(-4163)

$ws.Cells.Item(178, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(178, 2).Value = 'PASSED'
$ws.Cells.Item(178, 3).Value = 'firefox'
$ws.Cells.Item(178, 4).Formula = '="02.11.22"'
$ws.Cells.Item(178, 4).Copy()
$ws.Cells.Item(178, 4).PasteSpecial(-4163)

$ws.Cells.Item(179, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(179, 2).Value = 'PASSED'
$ws.Cells.Item(179, 3).Value = 'chrome'
$ws.Cells.Item(179, 4).Formula = '="02.11.22"'
$ws.Cells.Item(179, 4).Copy()
$ws.Cells.Item(179, 4).PasteSpecial(-4163)

$ws.Cells.Item(180, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(180, 2).Value = 'FAILED'
$ws.Cells.Item(180, 3).Value = 'firefox'
$ws.Cells.Item(180, 4).Formula = '="02.11.22"'
$ws.Cells.Item(180, 4).Copy()
$ws.Cells.Item(180, 4).PasteSpecial(-4163)

$ws.Cells.Item(181, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(181, 2).Value = 'PASSED'
$ws.Cells.Item(181, 3).Value = 'chrome'
$ws.Cells.Item(181, 4).Formula = '="02.11.22"'
$ws.Cells.Item(181, 4).Copy()
$ws.Cells.Item(181, 4).PasteSpecial(-4163)

$ws.Cells.Item(182, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(182, 2).Value = 'PASSED'
$ws.Cells.Item(182, 3).Value = 'firefox'
$ws.Cells.Item(182, 4).Formula = '="02.11.22"'
$ws.Cells.Item(182, 4).Copy()
$ws.Cells.Item(182, 4).PasteSpecial(-4163)

$ws.Cells.Item(183, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(183, 2).Value = 'PASSED'
$ws.Cells.Item(183, 3).Value = 'chrome'
$ws.Cells.Item(183, 4).Formula = '="02.11.22"'
$ws.Cells.Item(183, 4).Copy()
$ws.Cells.Item(183, 4).PasteSpecial(-4163)

$ws.Cells.Item(184, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(184, 2).Value = 'PASSED'
$ws.Cells.Item(184, 3).Value = 'firefox'
$ws.Cells.Item(184, 4).Formula = '="02.11.22"'
$ws.Cells.Item(184, 4).Copy()
$ws.Cells.Item(184, 4).PasteSpecial(-4163)

$ws.Cells.Item(185, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(185, 2).Value = 'FAILED'
$ws.Cells.Item(185, 3).Value = 'chrome'
$ws.Cells.Item(185, 4).Formula = '="02.11.22"'
$ws.Cells.Item(185, 4).Copy()
$ws.Cells.Item(185, 4).PasteSpecial(-4163)

$ws.Cells.Item(186, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(186, 2).Value = 'FAILED'
$ws.Cells.Item(186, 3).Value = 'firefox'
$ws.Cells.Item(186, 4).Formula = '="02.11.22"'
$ws.Cells.Item(186, 4).Copy()
$ws.Cells.Item(186, 4).PasteSpecial(-4163)

$ws.Cells.Item(187, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(187, 2).Value = 'FAILED'
$ws.Cells.Item(187, 3).Value = 'firefox'
$ws.Cells.Item(187, 4).Formula = '="02.11.22"'
$ws.Cells.Item(187, 4).Copy()
$ws.Cells.Item(187, 4).PasteSpecial(-4163)

$ws.Cells.Item(188, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(188, 2).Value = 'PASSED'
$ws.Cells.Item(188, 3).Value = 'chrome'
$ws.Cells.Item(188, 4).Formula = '="02.11.22"'
$ws.Cells.Item(188, 4).Copy()
$ws.Cells.Item(188, 4).PasteSpecial(-4163)

$ws.Cells.Item(189, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(189, 2).Value = 'PASSED'
$ws.Cells.Item(189, 3).Value = 'firefox'
$ws.Cells.Item(189, 4).Formula = '="02.11.22"'
$ws.Cells.Item(189, 4).Copy()
$ws.Cells.Item(189, 4).PasteSpecial(-4163)

$ws.Cells.Item(190, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(190, 2).Value = 'FAILED'
$ws.Cells.Item(190, 3).Value = 'chrome'
$ws.Cells.Item(190, 4).Formula = '="02.11.22"'
$ws.Cells.Item(190, 4).Copy()
$ws.Cells.Item(190, 4).PasteSpecial(-4163)

$ws.Cells.Item(191, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(191, 2).Value = 'PASSED'
$ws.Cells.Item(191, 3).Value = 'firefox'
$ws.Cells.Item(191, 4).Formula = '="02.11.22"'
$ws.Cells.Item(191, 4).Copy()
$ws.Cells.Item(191, 4).PasteSpecial(-4163)

$ws.Cells.Item(192, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(192, 2).Value = 'FAILED'
$ws.Cells.Item(192, 3).Value = 'chrome'
$ws.Cells.Item(192, 4).Formula = '="02.11.22"'
$ws.Cells.Item(192, 4).Copy()
$ws.Cells.Item(192, 4).PasteSpecial(-4163)

$ws.Cells.Item(193, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(193, 2).Value = 'PASSED'
$ws.Cells.Item(193, 3).Value = 'firefox'
$ws.Cells.Item(193, 4).Formula = '="02.11.22"'
$ws.Cells.Item(193, 4).Copy()
$ws.Cells.Item(193, 4).PasteSpecial(-4163)

$ws.Cells.Item(194, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(194, 2).Value = 'PASSED'
$ws.Cells.Item(194, 3).Value = 'chrome'
$ws.Cells.Item(194, 4).Formula = '="02.11.22"'
$ws.Cells.Item(194, 4).Copy()
$ws.Cells.Item(194, 4).PasteSpecial(-4163)

$ws.Cells.Item(195, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(195, 2).Value = 'FAILED'
$ws.Cells.Item(195, 3).Value = 'chrome'
$ws.Cells.Item(195, 4).Formula = '="02.11.22"'
$ws.Cells.Item(195, 4).Copy()
$ws.Cells.Item(195, 4).PasteSpecial(-4163)

$ws.Cells.Item(196, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(196, 2).Value = 'FAILED'
$ws.Cells.Item(196, 3).Value = 'firefox'
$ws.Cells.Item(196, 4).Formula = '="02.11.22"'
$ws.Cells.Item(196, 4).Copy()
$ws.Cells.Item(196, 4).PasteSpecial(-4163)

$ws.Cells.Item(197, 1).Value = 'Login with valid username and password'
$ws.Cells.Item(197, 2).Value = 'PASSED'
$ws.Cells.Item(197, 3).Value = 'chrome'
$ws.Cells.Item(197, 4).Formula = '="02.11.22"'
$ws.Cells.Item(197, 4).Copy()
$ws.Cells.Item(197, 4).PasteSpecial(-4163)

$ws.Cells.Item(198, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(198, 2).Value = 'PASSED'
$ws.Cells.Item(198, 3).Value = 'chrome'
$ws.Cells.Item(198, 4).Formula = '="02.11.22"'
$ws.Cells.Item(198, 4).Copy()
$ws.Cells.Item(198, 4).PasteSpecial(-4163)

$ws.Cells.Item(199, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(199, 2).Value = 'PASSED'
$ws.Cells.Item(199, 3).Value = 'chrome'
$ws.Cells.Item(199, 4).Formula = '="02.11.22"'
$ws.Cells.Item(199, 4).Copy()
$ws.Cells.Item(199, 4).PasteSpecial(-4163)

$ws.Cells.Item(200, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(200, 2).Value = 'PASSED'
$ws.Cells.Item(200, 3).Value = 'chrome'
$ws.Cells.Item(200, 4).Formula = '="02.11.22"'
$ws.Cells.Item(200, 4).Copy()
$ws.Cells.Item(200, 4).PasteSpecial(-4163)

$ws.Cells.Item(201, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(201, 2).Value = 'PASSED'
$ws.Cells.Item(201, 3).Value = 'chrome'
$ws.Cells.Item(201, 4).Formula = '="02.11.22"'
$ws.Cells.Item(201, 4).Copy()
$ws.Cells.Item(201, 4).PasteSpecial(-4163)

$ws.Cells.Item(202, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(202, 2).Value = 'PASSED'
$ws.Cells.Item(202, 3).Value = 'chrome'
$ws.Cells.Item(202, 4).Formula = '="02.11.22"'
$ws.Cells.Item(202, 4).Copy()
$ws.Cells.Item(202, 4).PasteSpecial(-4163)

$ws.Cells.Item(203, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(203, 2).Value = 'PASSED'
$ws.Cells.Item(203, 3).Value = 'chrome'
$ws.Cells.Item(203, 4).Formula = '="02.11.22"'
$ws.Cells.Item(203, 4).Copy()
$ws.Cells.Item(203, 4).PasteSpecial(-4163)

$ws.Cells.Item(204, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(204, 2).Value = 'PASSED'
$ws.Cells.Item(204, 3).Value = 'chrome'
$ws.Cells.Item(204, 4).Formula = '="02.11.22"'
$ws.Cells.Item(204, 4).Copy()
$ws.Cells.Item(204, 4).PasteSpecial(-4163)

$ws.Cells.Item(205, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(205, 2).Value = 'PASSED'
$ws.Cells.Item(205, 3).Value = 'chrome'
$ws.Cells.Item(205, 4).Formula = '="02.11.22"'
$ws.Cells.Item(205, 4).Copy()
$ws.Cells.Item(205, 4).PasteSpecial(-4163)

$ws.Cells.Item(206, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(206, 2).Value = 'FAILED'
$ws.Cells.Item(206, 3).Value = 'chrome'
$ws.Cells.Item(206, 4).Formula = '="02.11.22"'
$ws.Cells.Item(206, 4).Copy()
$ws.Cells.Item(206, 4).PasteSpecial(-4163)

$ws.Cells.Item(207, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(207, 2).Value = 'FAILED'
$ws.Cells.Item(207, 3).Value = 'chrome'
$ws.Cells.Item(207, 4).Formula = '="02.11.22"'
$ws.Cells.Item(207, 4).Copy()
$ws.Cells.Item(207, 4).PasteSpecial(-4163)

$ws.Cells.Item(208, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(208, 2).Value = 'PASSED'
$ws.Cells.Item(208, 3).Value = 'chrome'
$ws.Cells.Item(208, 4).Formula = '="02.11.22"'
$ws.Cells.Item(208, 4).Copy()
$ws.Cells.Item(208, 4).PasteSpecial(-4163)

$ws.Cells.Item(209, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(209, 2).Value = 'PASSED'
$ws.Cells.Item(209, 3).Value = 'chrome'
$ws.Cells.Item(209, 4).Formula = '="02.11.22"'
$ws.Cells.Item(209, 4).Copy()
$ws.Cells.Item(209, 4).PasteSpecial(-4163)

$ws.Cells.Item(210, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(210, 2).Value = 'FAILED'
$ws.Cells.Item(210, 3).Value = 'firefox'
$ws.Cells.Item(210, 4).Formula = '="02.11.22"'
$ws.Cells.Item(210, 4).Copy()
$ws.Cells.Item(210, 4).PasteSpecial(-4163)

$ws.Cells.Item(211, 1).Value = 'Login Failure with invalid username or password'
$ws.Cells.Item(211, 2).Value = 'FAILED'
$ws.Cells.Item(211, 3).Value = 'firefox'
$ws.Cells.Item(211, 4).Formula = '="02.11.22"'
$ws.Cells.Item(211, 4).Copy()
$ws.Cells.Item(211, 4).PasteSpecial(-4163)
